$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3: column I (praclen) 4 -> 5
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 5

# Update row 4 values and shift the old row 4 data down to row 5,
# inserting a new row 4 with the updated schedule.
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 61
$ws.Range("I4").Value = 5

# New row 5 (previously row 4's original data, now re-added as its own row)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 51
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = "train_dim1_2"

# Row 6 (previously row 5's data, shifted down, with praclen updated to 5)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 31
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim1_2"

$ws.Range("A3").Select()
